$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'24.519.11"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  -0.90%  '
$ws.Cells.Item(3, 4).Value = "'1.693.70"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  -0.41%  '
$ws.Cells.Item(4, 4).Value = "'1.001"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  -0.41%  '
$ws.Cells.Item(5, 4).Value = "'316.41"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -0.15%  '
$ws.Cells.Item(6, 4).Value = "'1.001"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -0.31%  '
$ws.Cells.Item(7, 4).Value = "'0.3909"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -0.88%  '
$ws.Cells.Item(8, 4).Value = "'0.4057"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +0.35%  '
$ws.Cells.Item(9, 4).Value = "'1.492"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -1.93%  '
$ws.Cells.Item(10, 4).Value = "'1.001"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -0.38%  '
$ws.Cells.Item(11, 4).Value = "'52.65"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -1.79%  '
$ws.Cells.Item(12, 4).Value = "'0.08785"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -1.01%  '
$ws.Cells.Item(13, 4).Value = "'26.68"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +12.58%  '
$ws.Cells.Item(14, 4).Value = "'7.532"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +0.77%  '
$ws.Cells.Item(15, 4).Value = "'8.136"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +0.08%  '
$ws.Cells.Item(16, 4).Value = "'0.00001351"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +1.92%  '
$ws.Cells.Item(17, 4).Value = "'1.686.93"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -1.10%  '
$ws.Cells.Item(18, 4).Value = "'98.10"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -1.65%  '
$ws.Cells.Item(19, 4).Value = "'0.07158"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +1.40%  '
$ws.Cells.Item(20, 4).Value = "'20.57"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +4.22%  '
$ws.Cells.Item(21, 4).Value = "'7.294"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +3.07%  '
$ws.Cells.Item(22, 4).Value = "'1.002"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -0.05%  '
$ws.Cells.Item(23, 4).Value = "'14.36"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -1.28%  '
$ws.Cells.Item(24, 4).Value = "'24.505.44"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -0.90%  '
$ws.Cells.Item(25, 4).Value = "'3.029"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -5.64%  '
$ws.Cells.Item(26, 4).Value = "'2.340"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -1.35%  '
$ws.Cells.Item(27, 4).Value = "'22.73"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -0.04%  '
$ws.Cells.Item(28, 4).Value = "'167.46"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +2.98%  '
$ws.Cells.Item(29, 4).Value = "'8.489"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -3.80%  '
$ws.Cells.Item(30, 4).Value = "'5.391"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +4.15%  '
$ws.Cells.Item(31, 4).Value = "'139.21"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +2.31%  '
$ws.Cells.Item(32, 4).Value = "'2.222"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +11.59%  '
$ws.Cells.Item(33, 4).Value = "'1.871.03"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -1.35%  '
$ws.Cells.Item(34, 4).Value = "'0.08788"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -1.76%  '
$ws.Cells.Item(35, 4).Value = "'7.324"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -5.92%  '
$ws.Cells.Item(36, 4).Value = "'1.041"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -3.46%  '
$ws.Cells.Item(37, 4).Value = "'0.02988"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +7.11%  '
$ws.Cells.Item(38, 4).Value = "'0.2787"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +1.08%  '
$ws.Cells.Item(39, 4).Value = "'10.96"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -1.03%  '
$ws.Cells.Item(40, 4).Value = "'0.09179"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +0.02%  '
$ws.Cells.Item(41, 4).Value = "'0.8079"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +4.57%  '
$ws.Cells.Item(42, 4).Value = "'14.19"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -2.58%  '
$ws.Cells.Item(43, 4).Value = "'1.476"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +0.87%  '
$ws.Cells.Item(44, 4).Value = "'17.67"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +10.00%  '
$ws.Cells.Item(45, 4).Value = "'2.686"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +4.48%  '
$ws.Cells.Item(46, 4).Value = "'0.7287"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +1.04%  '
$ws.Cells.Item(47, 4).Value = "'4.266"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +1.09%  '
$ws.Cells.Item(48, 4).Value = "'1.407"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +3.29%  '
$ws.Cells.Item(49, 4).Value = "'1.001"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -0.24%  '
$ws.Cells.Item(50, 4).Value = "'140.79"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -0.03%  '
$ws.Cells.Item(51, 4).Value = "'0.08188"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +2.42%  '
